$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 409.6270871618153

$ws.Range("B3").Value = 0.02258322285507441
$ws.Range("C3").Value = 0.05231270169004087
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 7.269056852830725
